$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 2743.074
$ws.Cells.Item(62, 9).Value = 3263.0557
$ws.Cells.Item(62, 10).Value = 1703.1111
$ws.Cells.Item(62, 11).Value = 3263.0557
$ws.Cells.Item(62, 12).Value = 1703.1111
$ws.Cells.Item(62, 13).Value = -2639.0557
$ws.Cells.Item(62, 14).Value = -2951.1111

$ws.Cells.Item(65, 8).Value = 2743.074
$ws.Cells.Item(65, 9).Value = 3263.0557
$ws.Cells.Item(65, 10).Value = 1703.1111
$ws.Cells.Item(65, 11).Value = 16315.2785
$ws.Cells.Item(65, 12).Value = 8515.5555
$ws.Cells.Item(65, 13).Value = -13195.2785
$ws.Cells.Item(65, 14).Value = -14755.5555

$ws.Cells.Item(129, 8).Value = 1810.0435
$ws.Cells.Item(129, 9).Value = 2897.75
$ws.Cells.Item(129, 10).Value = 1581.0526
$ws.Cells.Item(129, 11).Value = 8693.25
$ws.Cells.Item(129, 12).Value = 4743.1578
$ws.Cells.Item(129, 13).Value = -3693.25
$ws.Cells.Item(129, 14).Value = -14743.1578

$ws.Cells.Item(132, 8).Value = 38240.69
$ws.Cells.Item(132, 9).Value = 5451.1904
$ws.Cells.Item(132, 10).Value = 175956.6
$ws.Cells.Item(132, 11).Value = 16353.5712
$ws.Cells.Item(132, 12).Value = 527869.8
$ws.Cells.Item(132, 13).Value = -13823.5712
$ws.Cells.Item(132, 14).Value = -532929.8

$ws.Cells.Item(136, 8).Value = 36039
$ws.Cells.Item(136, 10).Value = 36039
$ws.Cells.Item(136, 12).Value = 36039
$ws.Cells.Item(136, 14).Value = -46239

$ws.Cells.Item(137, 8).Value = 5408.4375
$ws.Cells.Item(137, 9).Value = 1766.7333
$ws.Cells.Item(137, 10).Value = 8621.706
$ws.Cells.Item(137, 11).Value = 5300.199900000001
$ws.Cells.Item(137, 12).Value = 25865.118
$ws.Cells.Item(137, 13).Value = -2750.199900000001
$ws.Cells.Item(137, 14).Value = -30965.118

$ws.Cells.Item(138, 8).Value = 1324.9506
$ws.Cells.Item(138, 9).Value = 994.3582
$ws.Cells.Item(138, 10).Value = 2907.0715
$ws.Cells.Item(138, 11).Value = 2983.0746
$ws.Cells.Item(138, 12).Value = 8721.2145
$ws.Cells.Item(138, 13).Value = 2156.9254
$ws.Cells.Item(138, 14).Value = -19001.2145

$ws.Cells.Item(139, 8).Value = 41980
$ws.Cells.Item(139, 10).Value = 41980
$ws.Cells.Item(139, 12).Value = 41980
$ws.Cells.Item(139, 14).Value = -52260

$ws.Cells.Item(140, 8).Value = 50000
$ws.Cells.Item(140, 10).Value = 50000
$ws.Cells.Item(140, 12).Value = 50000
$ws.Cells.Item(140, 14).Value = -60360

$ws.Cells.Item(141, 8).Value = 1434.7869
$ws.Cells.Item(141, 9).Value = 964.0364
$ws.Cells.Item(141, 10).Value = 5750
$ws.Cells.Item(141, 11).Value = 2892.1092
$ws.Cells.Item(141, 12).Value = 17250
$ws.Cells.Item(141, 13).Value = 2287.8908
$ws.Cells.Item(141, 14).Value = -27610


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14819.207
$ws.Cells.Item(32, 9).Value = 11677.551
$ws.Cells.Item(32, 10).Value = 31923.777
$ws.Cells.Item(32, 11).Value = 11677.551
$ws.Cells.Item(32, 12).Value = 31923.777
$ws.Cells.Item(32, 13).Value = -11390.551
$ws.Cells.Item(32, 14).Value = -32497.777

$ws.Cells.Item(37, 8).Value = 38162.5
$ws.Cells.Item(37, 9).Value = 9000
$ws.Cells.Item(37, 10).Value = 47883.332
$ws.Cells.Item(37, 11).Value = 9000
$ws.Cells.Item(37, 12).Value = 47883.332
$ws.Cells.Item(37, 13).Value = -8727
$ws.Cells.Item(37, 14).Value = -48429.332

$ws.Cells.Item(41, 8).Value = 2454.6667
$ws.Cells.Item(41, 9).Value = 2454.6667
$ws.Cells.Item(41, 11).Value = 2454.6667
$ws.Cells.Item(41, 13).Value = -2040.6667

$ws.Cells.Item(61, 8).Value = 3018.6667
$ws.Cells.Item(61, 9).Value = 1958.0714
$ws.Cells.Item(61, 10).Value = 3693.5908
$ws.Cells.Item(61, 11).Value = 1958.0714
$ws.Cells.Item(61, 12).Value = 3693.5908
$ws.Cells.Item(61, 13).Value = -1746.0714
$ws.Cells.Item(61, 14).Value = -4117.5908

$ws.Cells.Item(107, 8).Value = 37324.332
$ws.Cells.Item(107, 10).Value = 37324.332
$ws.Cells.Item(107, 12).Value = 37324.332
$ws.Cells.Item(107, 14).Value = -45004.332

$ws.Cells.Item(136, 8).Value = 3018.6667
$ws.Cells.Item(136, 9).Value = 1958.0714
$ws.Cells.Item(136, 10).Value = 3693.5908
$ws.Cells.Item(136, 11).Value = 5874.2142
$ws.Cells.Item(136, 12).Value = 11080.7724
$ws.Cells.Item(136, 13).Value = -3324.2142
$ws.Cells.Item(136, 14).Value = -16180.7724


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1720.8096
$ws.Cells.Item(107, 9).Value = 1599.4
$ws.Cells.Item(107, 10).Value = 2024.3334
$ws.Cells.Item(107, 11).Value = 1599.4
$ws.Cells.Item(107, 12).Value = 2024.3334
$ws.Cells.Item(107, 13).Value = 320.5999999999999
$ws.Cells.Item(107, 14).Value = -5864.3334

$ws.Cells.Item(137, 8).Value = 33136
$ws.Cells.Item(137, 10).Value = 33136
$ws.Cells.Item(137, 12).Value = 33136
$ws.Cells.Item(137, 14).Value = -43336


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 167749.61
$ws.Cells.Item(31, 9).Value = 2273.9092
$ws.Cells.Item(31, 10).Value = 192019.38
$ws.Cells.Item(31, 11).Value = 2273.9092
$ws.Cells.Item(31, 12).Value = 192019.38
$ws.Cells.Item(31, 13).Value = -1978.9092
$ws.Cells.Item(31, 14).Value = -192609.38

$ws.Cells.Item(34, 8).Value = 167749.61
$ws.Cells.Item(34, 9).Value = 2273.9092
$ws.Cells.Item(34, 10).Value = 192019.38
$ws.Cells.Item(34, 11).Value = 2273.9092
$ws.Cells.Item(34, 12).Value = 192019.38
$ws.Cells.Item(34, 13).Value = -2071.9092
$ws.Cells.Item(34, 14).Value = -192423.38

$ws.Cells.Item(58, 8).Value = 1817.3214
$ws.Cells.Item(58, 9).Value = 1515.5834
$ws.Cells.Item(58, 10).Value = 2043.625
$ws.Cells.Item(58, 11).Value = 1515.5834
$ws.Cells.Item(58, 12).Value = 2043.625
$ws.Cells.Item(58, 13).Value = -1312.5834
$ws.Cells.Item(58, 14).Value = -2449.625

$ws.Cells.Item(107, 8).Value = 409.5091
$ws.Cells.Item(107, 9).Value = 352.78125
$ws.Cells.Item(107, 10).Value = 488.43478
$ws.Cells.Item(107, 11).Value = 352.78125
$ws.Cells.Item(107, 12).Value = 488.43478
$ws.Cells.Item(107, 13).Value = 1567.21875
$ws.Cells.Item(107, 14).Value = -4328.43478

$ws.Cells.Item(115, 8).Value = 28236.75
$ws.Cells.Item(115, 10).Value = 28236.75
$ws.Cells.Item(115, 12).Value = 28236.75
$ws.Cells.Item(115, 14).Value = -30586.75

$ws.Cells.Item(132, 8).Value = 76226.63
$ws.Cells.Item(132, 9).Value = 2054.7693
$ws.Cells.Item(132, 10).Value = 236932.33
$ws.Cells.Item(132, 11).Value = 6164.3079
$ws.Cells.Item(132, 12).Value = 710796.99
$ws.Cells.Item(132, 13).Value = -3634.3079
$ws.Cells.Item(132, 14).Value = -715856.99

$ws.Cells.Item(134, 8).Value = 319702.78
$ws.Cells.Item(134, 9).Value = 1142.0741
$ws.Cells.Item(134, 10).Value = 825652.1
$ws.Cells.Item(134, 11).Value = 3426.2223
$ws.Cells.Item(134, 12).Value = 2476956.3
$ws.Cells.Item(134, 13).Value = -891.2223000000004
$ws.Cells.Item(134, 14).Value = -2482026.3

$ws.Cells.Item(136, 8).Value = 1817.3214
$ws.Cells.Item(136, 9).Value = 1515.5834
$ws.Cells.Item(136, 10).Value = 2043.625
$ws.Cells.Item(136, 11).Value = 4546.7502
$ws.Cells.Item(136, 12).Value = 6130.875
$ws.Cells.Item(136, 13).Value = -1996.7502
$ws.Cells.Item(136, 14).Value = -11230.875


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1807.4117
$ws.Cells.Item(113, 9).Value = 1922.2222
$ws.Cells.Item(113, 10).Value = 1678.25
$ws.Cells.Item(113, 11).Value = 1922.2222
$ws.Cells.Item(113, 12).Value = 1678.25
$ws.Cells.Item(113, 13).Value = 247.7778000000001
$ws.Cells.Item(113, 14).Value = -6018.25

$ws.Cells.Item(136, 8).Value = 33150.285
$ws.Cells.Item(136, 10).Value = 33150.285
$ws.Cells.Item(136, 12).Value = 99450.85500000001
$ws.Cells.Item(136, 14).Value = -104550.855


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(45, 8).Value = 9996.666999999999
$ws.Cells.Item(45, 9).Value = 10000
$ws.Cells.Item(45, 10).Value = 9995
$ws.Cells.Item(45, 11).Value = 10000
$ws.Cells.Item(45, 12).Value = 9995
$ws.Cells.Item(45, 13).Value = -9593
$ws.Cells.Item(45, 14).Value = -10809

$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 13).ClearContents()

$ws.Cells.Item(75, 8).Value = 37173
$ws.Cells.Item(75, 10).Value = 37173
$ws.Cells.Item(75, 12).Value = 37173
$ws.Cells.Item(75, 14).Value = -39045

$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 13).ClearContents()

$ws.Cells.Item(78, 8).Value = 37173
$ws.Cells.Item(78, 10).Value = 37173
$ws.Cells.Item(78, 12).Value = 111519
$ws.Cells.Item(78, 14).Value = -120879

$ws.Cells.Item(132, 8).Value = 2287.7415
$ws.Cells.Item(132, 9).Value = 1524.5366
$ws.Cells.Item(132, 11).Value = 4573.6098
$ws.Cells.Item(132, 13).Value = -2043.6098

$ws.Cells.Item(136, 8).Value = 1724.5
$ws.Cells.Item(136, 9).Value = 1144
$ws.Cells.Item(136, 10).Value = 2930.1538
$ws.Cells.Item(136, 11).Value = 3432
$ws.Cells.Item(136, 12).Value = 8790.4614
$ws.Cells.Item(136, 13).Value = -882
$ws.Cells.Item(136, 14).Value = -13890.4614


$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 250034670
$ws.Cells.Item(119, 10).Value = 250034670
$ws.Cells.Item(119, 12).Value = 250034670
$ws.Cells.Item(119, 14).Value = -250044346

$ws.Cells.Item(135, 8).Value = 33379.8
$ws.Cells.Item(135, 10).Value = 33379.8
$ws.Cells.Item(135, 12).Value = 33379.8
$ws.Cells.Item(135, 14).Value = -43519.8

$ws.Cells.Item(136, 8).Value = 22278.617
$ws.Cells.Item(136, 9).Value = 29203.258
$ws.Cells.Item(136, 10).Value = 2081.75
$ws.Cells.Item(136, 11).Value = 87609.774
$ws.Cells.Item(136, 12).Value = 6245.25
$ws.Cells.Item(136, 13).Value = -85059.774
$ws.Cells.Item(136, 14).Value = -11345.25

